$wb = $excel.ActiveWorkbook

# OFF sheet - Week 16 logged stats (row 2 = Home)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 146
$wsOff.Range("C2").Value = 93
$wsOff.Range("D2").Value = 47
$wsOff.Range("E2").Value = 23

# DEF sheet - Week 16 logged stats (row 2 = Home)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 215
$wsDef.Range("C2").Value = 147
$wsDef.Range("D2").Value = 35
$wsDef.Range("F2").Value = 4
$wsDef.Range("G2").Value = 3
